$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its values as literal text, matching the
# source data which stores prices as formatted strings (e.g. thousands separated
# with dots) rather than numbers. Without this, Excel would auto-convert numeric-
# looking strings (like "1.000" or "41.97") into actual numbers and drop formatting.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '26.895.71'
$ws.Range("E2").Value = '  -1.01%  '

# Row 3
$ws.Range("D3").Value = '1.739.87'
$ws.Range("E3").Value = '  +1.24%  '

# Row 4
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").Value = '310.70'
$ws.Range("E5").Value = '  +0.24%  '

# Row 6
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.03%  '

# Row 7
$ws.Range("D7").Value = '0.5021'
$ws.Range("E7").Value = '  +7.43%  '

# Row 8
$ws.Range("D8").Value = '0.3574'
$ws.Range("E8").Value = '  +4.08%  '

# Row 9
$ws.Range("D9").Value = '41.97'
$ws.Range("E9").Value = '  -0.24%  '

# Row 10
$ws.Range("D10").Value = '0.07241'
$ws.Range("E10").Value = '  -0.25%  '

# Row 11
$ws.Range("D11").Value = '1.059'
$ws.Range("E11").Value = '  +1.26%  '

# Row 12
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  -0.02%  '

# Row 13
$ws.Range("D13").Value = '20.14'
$ws.Range("E13").Value = '  +1.21%  '

# Row 14
$ws.Range("D14").Value = '5.948'
$ws.Range("E14").Value = '  +1.65%  '

# Row 15
$ws.Range("D15").Value = '1.738.19'
$ws.Range("E15").Value = '  +1.56%  '

# Row 16
$ws.Range("D16").Value = '6.805'
$ws.Range("E16").Value = '  -1.04%  '

# Row 17
$ws.Range("D17").Value = '86.52'
$ws.Range("E17").Value = '  -2.60%  '

# Row 18
$ws.Range("D18").Value = '0.00001032'
$ws.Range("E18").Value = '  -0.59%  '

# Row 19
$ws.Range("E19").Value = '  +1.01%  '

# Row 20
$ws.Range("E20").Value = '  -0.02%  '

# Row 21
$ws.Range("D21").Value = '16.52'
$ws.Range("E21").Value = '  -0.24%  '

# Row 22
$ws.Range("D22").Value = '5.738'
$ws.Range("E22").Value = '  +1.65%  '

# Row 23
$ws.Range("D23").Value = '26.964.47'
$ws.Range("E23").Value = '  -0.82%  '

# Row 24
$ws.Range("E24").Value = '  +3.40%  '

# Row 25
$ws.Range("D25").Value = '2.043'
$ws.Range("E25").Value = '  -4.25%  '

# Row 26
$ws.Range("D26").Value = '152.80'
$ws.Range("E26").Value = '  -2.24%  '

# Row 27
$ws.Range("D27").Value = '19.81'
$ws.Range("E27").Value = '  +2.07%  '

# Row 28
$ws.Range("D28").Value = '1.942.86'
$ws.Range("E28").Value = '  +1.69%  '

# Row 29
$ws.Range("D29").Value = '2.208'
$ws.Range("E29").Value = '  +4.12%  '

# Row 30
$ws.Range("D30").Value = '119.62'
$ws.Range("E30").Value = '  -0.22%  '

# Row 31
$ws.Range("D31").Value = '1.040'
$ws.Range("E31").Value = '  +1.63%  '

# Row 32
$ws.Range("D32").Value = '0.09550'
$ws.Range("E32").Value = '  +4.38%  '

# Row 33
$ws.Range("D33").Value = '3.579'
$ws.Range("E33").Value = '  -0.38%  '

# Row 34
$ws.Range("D34").Value = '5.362'
$ws.Range("E34").Value = '  +0.53%  '

# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.05892'
$ws.Range("E35").Value = '  +1.01%  '

# Row 36
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.02175'
$ws.Range("E36").Value = '  -0.88%  '

# Row 37
$ws.Range("D37").Value = '10.99'
$ws.Range("E37").Value = '  -0.56%  '

# Row 38
$ws.Range("D38").Value = '0.2000'
$ws.Range("E38").Value = '  +0.09%  '

# Row 39
$ws.Range("D39").Value = '1.423'
$ws.Range("E39").Value = '  +2.44%  '

# Row 40
$ws.Range("D40").Value = '4.758'
$ws.Range("E40").Value = '  +0.53%  '

# Row 41
$ws.Range("D41").Value = '0.6027'
$ws.Range("E41").Value = '  +1.70%  '

# Row 42
$ws.Range("D42").Value = '1.111'
$ws.Range("E42").Value = '  -1.48%  '

# Row 43
$ws.Range("D43").Value = '7.572'
$ws.Range("E43").Value = '  +1.28%  '

# Row 44
$ws.Range("D44").Value = '12.90'
$ws.Range("E44").Value = '  +2.40%  '

# Row 45
$ws.Range("E45").Value = '  +0.67%  '

# Row 46
$ws.Range("D46").Value = '0.5643'
$ws.Range("E46").Value = '  -0.27%  '

# Row 47
$ws.Range("D47").Value = '119.97'
$ws.Range("E47").Value = '  +0.74%  '

# Row 48
$ws.Range("D48").Value = '1.841'
$ws.Range("E48").Value = '  -0.24%  '

# Row 49
$ws.Range("D49").Value = '1.098'
$ws.Range("E49").Value = '  +1.19%  '

# Row 50
$ws.Range("D50").Value = '0.06649'
$ws.Range("E50").Value = '  -0.21%  '

# Row 51
$ws.Range("D51").Value = '1.000'
$ws.Range("E51").Value = '  +0.02%  '
